# Methodos numericos.xlsx -- add Newton-Raphson / Secante sheet, plus misc
# selection / formatting touch-ups that a user would have made while at it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Hoja1: the percent-error column I3 switches from the custom
#    "0.0000%" look to the builtin "0.00%" style, and the selection
#    that was left behind moves to E17.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("I3").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 2. Hoja2: no longer the tab that was active when the file was saved,
#    and the lingering selection moves back to B3.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("B3").Select()

# ---------------------------------------------------------------------
# 3. Three new sheets: Hoja3, Hoja4 (left blank) and Hoja5, which gets
#    the Newton-Raphson / Secante method tables.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Hoja3"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Hoja4"

$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "Hoja5"

# --- Hoja5 content -----------------------------------------------------

# Titles
$ws5.Range("A1:G1").HorizontalAlignment = -4108
$ws5.Range("A1:F1").Merge()
$ws5.Range("A1").Value = "Metodo  Newton-Rahpson"

$ws5.Range("I1:R1").HorizontalAlignment = -4108
$ws5.Range("I1:M1").Merge()
$ws5.Range("I1").Value = "Metodo Secante"

# Header row (row 2)
$ws5.Range("A2:F2").HorizontalAlignment = -4108
$ws5.Range("A2").Value = "I"
$ws5.Range("B2").Value = "Xi"
$ws5.Range("C2").Value = "F(X)"
$ws5.Range("D2").Value = "F'(X)"
$ws5.Range("E2").Value = "F(X)/F'(X)"
$ws5.Range("F2").Value = "ERRORa%"

$ws5.Range("I2:M2").HorizontalAlignment = -4108
$ws5.Range("I2").Value = "Xi"
$ws5.Range("J2").Value = "F(Xi)"
$ws5.Range("K2").Value = "Xi-1"
$ws5.Range("L2").Value = "F(Xi-1)"
$ws5.Range("M2").Value = "ERRORa%"

# Newton-Raphson iterations
$ws5.Range("A3").Value = 1
$ws5.Range("B3").Value = 0
$ws5.Range("C3").Formula = "=(EXP(-B3))-B3"
$ws5.Range("D3").Formula = "=EXP(-B3)+1"
$ws5.Range("E3").Formula = "=C3/D3"

$ws5.Range("A4").Value = 2
$ws5.Range("B4").Formula = "=B3-(C3/D3)"
$ws5.Range("C4").Formula = "=(EXP(-B4))-B4"
$ws5.Range("D4").Formula = "=-(EXP(-B4)+1)"
$ws5.Range("E4").Formula = "=C4/D4"
$ws5.Range("F4").Formula = "=ABS((B4-B3)/B4)"

$ws5.Range("A5").Value = 3
$ws5.Range("B5").Formula = "=IF(F4<1,B4-(C4/D4),0)"
$ws5.Range("C5").Formula = "=(EXP(-B5))-B5"
$ws5.Range("D5").Formula = "=-(EXP(-B5)+1)"
$ws5.Range("E5").Formula = "=C5/D5"
$ws5.Range("F5").Formula = "=ABS((B5-B4)/B5)"

$ws5.Range("A6").Value = 4
$ws5.Range("B6").Formula = "=IF(F5<1,B5-(C5/D5),0)"
$ws5.Range("C6").Formula = "=(EXP(-B6))-B6"
$ws5.Range("D6").Formula = "=-(EXP(-B6)+1)"
$ws5.Range("E6").Formula = "=C6/D6"
$ws5.Range("F6").Formula = "=ABS((B6-B5)/B6)"

$ws5.Range("A7").Value = 5
$ws5.Range("A8").Value = 6
$ws5.Range("A9").Value = 7
$ws5.Range("A10").Value = 8
$ws5.Range("A11").Value = 9
$ws5.Range("A12").Value = 10

# Number formats / alignment for the worked columns, matching the rest
# of the workbook's look (2-decimal, centred; percent for the error col).
$ws5.Range("B4:B14").NumberFormat = "0.00"
$ws5.Range("B4:B14").HorizontalAlignment = -4108
$ws5.Range("C10:E14").NumberFormat = "0.00"
$ws5.Range("C10:E14").HorizontalAlignment = -4108
$ws5.Range("F4:F10").NumberFormat = "0.00%"
$ws5.Range("F11:F14").NumberFormat = "0.00%"

$ws5.Range("I4:L14").NumberFormat = "0.00"
$ws5.Range("I4:L14").HorizontalAlignment = -4108
$ws5.Range("M4:M14").NumberFormat = "0.00%"
$ws5.Range("I10").NumberFormat = "0.00"
$ws5.Range("I10").HorizontalAlignment = -4108

$ws5.Range("B4").Select()
$ws5.Activate()
